# Generate Report for Handback
# Populates the "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" columns on the per-locale sheets now that the handback round-trip
# completed, and flips the Status column from "Ready for handoff" to
# "Handed back: in sync with en-US" everywhere it's shown (Overview + each
# locale sheet).

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

$mdDisplay   = "958c71ad-ad56-40fa-bf50-211a36cca7b3.md"
$mdUrl       = "https://github.com/OpenLocalizationTest/oltest/blob/1559a4d5828084ad6a644c22f5b305000986c580/e2e/958c71ad-ad56-40fa-bf50-211a36cca7b3.md"

$zhXlfDisplay = "958c71ad-ad56-40fa-bf50-211a36cca7b3.d3dd10a1426c8998aa01023e28b816c216feb709.zh-cn.xlf"
$zhXlfUrl     = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0a2ee87cead2bfc59c84454470ea0330ddb8cfb8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/958c71ad-ad56-40fa-bf50-211a36cca7b3.d3dd10a1426c8998aa01023e28b816c216feb709.zh-cn.xlf"

$deXlfDisplay = "958c71ad-ad56-40fa-bf50-211a36cca7b3.d3dd10a1426c8998aa01023e28b816c216feb709.de-de.xlf"
$deXlfUrl     = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/da8a87eb20590ce4a45635768ef22fd02eef1c26/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/958c71ad-ad56-40fa-bf50-211a36cca7b3.d3dd10a1426c8998aa01023e28b816c216feb709.de-de.xlf"

$zhHandbackTime = "2016-01-27 08:41:14"
$deHandbackTime = "2016-01-27 08:41:34"

# --- Overview sheet: flip the status text shown for each locale/file row ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusHandedBack
$wsOverview.Range("C2").Value = $statusHandedBack
$wsOverview.Range("B3").Value = $statusHandedBack
$wsOverview.Range("C3").Value = $statusHandedBack

function Set-HandbackColumns {
    param(
        $ws,
        $fileDisplay,
        $fileUrl,
        $handbackTime
    )

    # Status -> handed back, in sync with en-US
    $ws.Range("B2").Value = $statusHandedBack
    $ws.Range("B3").Value = $statusHandedBack

    # Latest Target File (E) / Latest Handback File (F) are populated with the
    # same source-md / xlf links used for the handoff, now that they came back.
    $ws.Hyperlinks.Add($ws.Range("E2"), $mdUrl, "", "", $mdDisplay)
    $ws.Hyperlinks.Add($ws.Range("F2"), $fileUrl, "", "", $fileDisplay)
    $ws.Hyperlinks.Add($ws.Range("E3"), $mdUrl, "", "", $mdDisplay)
    $ws.Hyperlinks.Add($ws.Range("F3"), $fileUrl, "", "", $fileDisplay)

    # Match the workbook's existing hyperlink look (underlined, cornflower blue)
    # rather than the theme-default hyperlink color.
    $ws.Range("E2").Font.Underline = $true
    $ws.Range("E2").Font.Color = 0xED9564
    $ws.Range("F2").Font.Underline = $true
    $ws.Range("F2").Font.Color = 0xED9564
    $ws.Range("E3").Font.Underline = $true
    $ws.Range("E3").Font.Color = 0xED9564
    $ws.Range("F3").Font.Underline = $true
    $ws.Range("F3").Font.Color = 0xED9564

    # Latest Handback DateTime (G)
    $ws.Range("G2").Value = $handbackTime
    $ws.Range("G3").Value = $handbackTime
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Set-HandbackColumns $wsZhCn $zhXlfDisplay $zhXlfUrl $zhHandbackTime

$wsDeDe = $wb.Worksheets.Item("de-de")
Set-HandbackColumns $wsDeDe $deXlfDisplay $deXlfUrl $deHandbackTime
